$wb = $excel.ActiveWorkbook

# --- SignInData sheet: clear row 3's data (second sign-in record) -------
# Row 3 (nivakalita744@gmail.com / kalitaniva@1234 / Welcome) is removed:
# A3/B3 become blank (keeping their existing "Hyperlink" cell style) and
# C3 is dropped entirely. The two hyperlinks that lived on A3/B3 must go
# too, but their mailto addresses for A2/B2 must be preserved.
$signIn = $wb.Worksheets.Item("SignInData")

# This runtime's Hyperlink.Delete only takes effect when the whole
# collection is cleared, so capture the two links we want to KEEP, wipe
# everything, then recreate just those two (re-applying their original
# "Hyperlink" cell style afterwards since Add() re-stamps it anyway).
$signIn.Hyperlinks.Delete()
$signIn.Hyperlinks.Add($signIn.Range("A2"), "mailto:nayangogoi744@gmail.com")
$signIn.Hyperlinks.Add($signIn.Range("B2"), "mailto:SAG*@dem01234")
$signIn.Range("A2:B2").Style = "Hyperlink"

# Now blank out row 3 (content only, formatting of A3/B3 stays; C3 had no
# special style so it disappears completely once its value is cleared).
$signIn.Range("A3:C3").ClearContents()

# The row that was edited is now the selection, and this sheet becomes the
# active/selected tab of the workbook (previously it was AccountData).
$signIn.Activate()
$signIn.Rows.Item(3).Select()
